$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the label in B8 from "id" to "insuredId"
$ws.Range("B8").Value = "insuredId"

# Update the current selection to B9 (as last interacted cell)
$ws.Range("B9").Select()
